$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "320.93"
Set-TextValue $ws "E2" "8.03%"
Set-TextValue $ws "D3" "48.00"
Set-TextValue $ws "E3" "14.28%"
Set-TextValue $ws "D4" "5.259"
Set-TextValue $ws "E4" "4.84%"
Set-TextValue $ws "D5" "0.08093"
Set-TextValue $ws "E5" "7.48%"
Set-TextValue $ws "D6" "4.569"
Set-TextValue $ws "E6" "4.16%"
Set-TextValue $ws "D7" "1.670"
Set-TextValue $ws "E7" "6.18%"
Set-TextValue $ws "D8" "1.093"
Set-TextValue $ws "E8" "17.68%"
Set-TextValue $ws "D9" "0.1322"
Set-TextValue $ws "E9" "11.33%"
Set-TextValue $ws "D10" "0.1957"
Set-TextValue $ws "E10" "6.93%"
Set-TextValue $ws "D11" "0.09588"
Set-TextValue $ws "E11" "7.89%"
Set-TextValue $ws "D12" "0.04532"
Set-TextValue $ws "E12" "10.92%"
Set-TextValue $ws "E13" "0.16%"
Set-TextValue $ws "D14" "0.001337"
Set-TextValue $ws "E14" "3.47%"
Set-TextValue $ws "D15" "0.005799"
Set-TextValue $ws "E15" "-3.07%"
Set-TextValue $ws "D16" "3.381"
Set-TextValue $ws "E16" "0.66%"
Set-TextValue $ws "D17" "2.428"
Set-TextValue $ws "E17" "1.10%"
Set-TextValue $ws "E18" "2.33%"
Set-TextValue $ws "D19" "8.180"
Set-TextValue $ws "E19" "1.63%"
Set-TextValue $ws "D20" "0.1397"
Set-TextValue $ws "E20" "-0.98%"
Set-TextValue $ws "D21" "0.2926"
Set-TextValue $ws "E21" "-11.40%"
Set-TextValue $ws "D22" "0.04316"
Set-TextValue $ws "E22" "4.72%"
Set-TextValue $ws "D23" "0.001306"
Set-TextValue $ws "E23" "3.26%"
Set-TextValue $ws "E24" "9.43%"
Set-TextValue $ws "D25" "0.0001349"
Set-TextValue $ws "E25" "9.63%"
Set-TextValue $ws "D38" "0.02776"
Set-TextValue $ws "E38" "15.45%"
Set-TextValue $ws "D39" "0.05541"
Set-TextValue $ws "E39" "6.50%"
Set-TextValue $ws "D40" "0.006286"
Set-TextValue $ws "E40" "-0.27%"
Set-TextValue $ws "D41" "0.007782"
Set-TextValue $ws "E41" "-0.51%"
Set-TextValue $ws "D42" "0.1445"
Set-TextValue $ws "E42" "8.82%"
Set-TextValue $ws "E43" "3.55%"
Set-TextValue $ws "D44" "0.008834"
Set-TextValue $ws "E44" "26.32%"
Set-TextValue $ws "D45" "0.3512"
Set-TextValue $ws "E45" "18.72%"
Set-TextValue $ws "D46" "0.00006845"
Set-TextValue $ws "E46" "5.93%"
Set-TextValue $ws "E47" "-0.27%"
Set-TextValue $ws "D48" "0.05330"
Set-TextValue $ws "E48" "2.10%"
Set-TextValue $ws "E49" "-5.06%"
Set-TextValue $ws "D50" "0.00002095"
Set-TextValue $ws "E50" "-0.27%"
Set-TextValue $ws "E51" "-0.27%"
